# BALP 1.1.1 and history file updates
# Update the "Metadata" property sheet of the UserAgentTypes CodeSystem workbook:
#   - Version bumps from 1.1.0 to 1.1.1
#   - Experimental now explicitly shows "false"
#   - Date moves to the new publication date
#   - Case Sensitive now explicitly shows "false"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "1.1.1"

# Leading apostrophe forces these to be stored as literal text "false"
# instead of being auto-converted to a Boolean TRUE/FALSE value.
$ws.Range("B7").Value = "'false"
$ws.Range("B8").Value = "2022-10-21T09:04:31-05:00"
$ws.Range("B17").Value = "'false"
